$d = $word.ActiveDocument

# "Division: D20A" -> "Division: D20B"
$d.Content.Find.Execute("Division: D20A", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Division: D20B", 2)

# "Roll no: 64" -> "Roll no: 60"
$d.Content.Find.Execute("Roll no: 64", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Roll no: 60", 2)
